# Duplicate row 2 (the sample "jhi_user" record) into row 3 on the
# "jhi_user" sheet, extending the used range from A1:N2 to A1:N3.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jhi_user")

$ws.Range("B2:N2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial() | Out-Null
